# Append " (Changed main)" to the end of the first paragraph's text, right
# after "This is a Microsoft word document." — turning it into:
#   "This is a Microsoft word document. (Changed main)"

$d = $word.ActiveDocument
$firstPara = $d.Paragraphs.First

if ($firstPara.Range.Text -notlike "*Changed main*") {
    $d.Content.Find.Execute(
        "This is a Microsoft word document.", $true, $false, $false, $false, $false,
        $true, 1, $false,
        "This is a Microsoft word document. (Changed main)", 2
    ) | Out-Null
}

# Fallback / verification path: if, for whatever reason, Find & Replace did
# not perform the substitution, locate the first paragraph directly and
# append the pieces so the end result is still correct.
$firstPara = $d.Paragraphs.First
if ($firstPara.Range.Text -notlike "*Changed main*") {
    $endPos = $firstPara.Range.End - 1

    $r1 = $d.Range($endPos, $endPos)
    $r1.InsertAfter(" (")

    $endPos = $firstPara.Range.End - 1
    $r2 = $d.Range($endPos, $endPos)
    $r2.InsertAfter("Changed main")

    $endPos = $firstPara.Range.End - 1
    $r3 = $d.Range($endPos, $endPos)
    $r3.InsertAfter(")")
}

Write-Output $d.Paragraphs.First.Range.Text
